# Translations.xlsx - add FULL translations for the "Selected Obj Panel"
#
# Strategy (derived from the shared-strings append order baked into the
# diff): the author first appended the 26 new attribute rows at the
# bottom of the sheet (rows 65-90, columns A/B/C, in that left-to-right
# order - except row 81 where the Spanish/English/Id triplet was typed
# C,B,A), and only afterwards went back and inserted a single new row
# right under the header block (row 63) for the "Collision" entry. We
# reproduce writes in that exact chronological order so new entries land
# in the shared string table with the same indices as the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Make room for the new "Collision" row right after the header block
#    (row 63) FIRST, so it pushes only the pre-existing row 63
#    (tooltip.SetActiveAtStartToggle) down to row 64 and does not
#    disturb the brand-new rows we are about to append below it. The
#    inserted row's cells are left blank for now - they are filled in
#    last (step 2) so the new "Collision" string lands at the end of
#    the shared-string table, matching the saved file.
# ---------------------------------------------------------------------

$ws.Rows(63).Insert()

# ---------------------------------------------------------------------
# 1. Append the new attribute rows at the bottom (rows 65-90).
# ---------------------------------------------------------------------

$ws.Range("A65").Value = "ColorHex"
$ws.Range("B65").Value = "Color (Hex)"
$ws.Range("C65").Value = "Color (Hex)"

$ws.Range("A66").Value = "Intensity"
$ws.Range("B66").Value = "Intensity"
$ws.Range("C66").Value = "Intensidad"

$ws.Range("A67").Value = "Range"
$ws.Range("B67").Value = "Range"
$ws.Range("C67").Value = "Rango"

$ws.Range("A68").Value = "ActivateOnStart"
$ws.Range("B68").Value = "Activate On Start"
$ws.Range("C68").Value = "Activar Al Inicio"

$ws.Range("A69").Value = "Damage"
$ws.Range("B69").Value = "Damage"
$ws.Range("C69").Value = "Daño"

$ws.Range("A70").Value = "TravelBack"
$ws.Range("B70").Value = "Travel Back"
$ws.Range("C70").Value = "Regresarse"

$ws.Range("A71").Value = "AddWaypoint"
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").VerticalAlignment = -4160
$ws.Range("B71").WrapText = $true
$ws.Range("B71").Value = "+ Add Waypoint"
$ws.Range("C71").NumberFormat = "@"
$ws.Range("C71").VerticalAlignment = -4160
$ws.Range("C71").WrapText = $true
$ws.Range("C71").Value = "+ Añadir Waypoint"

$ws.Range("A72").Value = "WaitTime"
$ws.Range("B72").Value = "Wait Time"
$ws.Range("C72").Value = "Tiempo de Espera"

$ws.Range("A73").Value = "InitialState"
$ws.Range("B73").Value = "Initial State"
$ws.Range("C73").Value = "Estado Inicial"

$ws.Range("A74").Value = "DEACTIVATED"
$ws.Range("B74").Value = "DEACTIVATED"
$ws.Range("C74").Value = "DESACTIVADO"

$ws.Range("A75").Value = "ACTIVATED"
$ws.Range("B75").Value = "ACTIVATED"
$ws.Range("C75").Value = "ACTIVADO"

$ws.Range("A76").Value = "UNUSABLE"
$ws.Range("B76").Value = "UNUSABLE"
$ws.Range("C76").Value = "UNUSABLE"

$ws.Range("A77").Value = "UsableOnce"
$ws.Range("B77").Value = "Usable Once"
$ws.Range("C77").Value = "Usable Una Vez"

$ws.Range("A78").Value = "CanBeShotByTaser"
$ws.Range("B78").Value = "Can be shot by Taser"
$ws.Range("C78").Value = "Puede ser activado por el Taser"

$ws.Range("A79").Value = "ManageEvents"
$ws.Range("B79").Value = "Manage Events"
$ws.Range("C79").Value = "Administrar Eventos"

$ws.Range("A80").Value = "RespawnTime"
$ws.Range("B80").Value = "Respawn Time"
$ws.Range("C80").Value = "Tiempo de Reaparición"

# Row 81 was authored Spanish -> English -> Id, unlike the rest, so the
# shared-string indices for this row land out of the usual order.
$ws.Range("C81").Value = "Muerte Instantanea"
$ws.Range("B81").Value = "Instant Kill"
$ws.Range("A81").Value = "InstantKill"

$ws.Range("A82").Value = "Constant"
$ws.Range("B82").Value = "Constant"
$ws.Range("C82").Value = "Constante"

$ws.Range("A83").Value = "OnlyOnce"
$ws.Range("B83").Value = "Only Once"
$ws.Range("C83").Value = "Una sola Vez"

$ws.Range("A84").Value = "ScreenColor"
$ws.Range("B84").Value = "Screen Color"
$ws.Range("C84").Value = "Color de la Pantalla"

$ws.Range("A85").Value = "CYAN"
$ws.Range("B85").Value = "CYAN"
$ws.Range("C85").Value = "CYAN"

$ws.Range("A86").Value = "GREEN"
$ws.Range("B86").Value = "GREEN"
$ws.Range("C86").Value = "VERDE"

$ws.Range("A87").Value = "RED"
$ws.Range("B87").Value = "RED"
$ws.Range("C87").Value = "ROJO"

$ws.Range("A88").Value = "InvisibleMesh"
$ws.Range("B88").Value = "Invisible Mesh"
$ws.Range("C88").Value = "Malla Invisible"

$ws.Range("A89").Value = "InvertTextWithGravity"
$ws.Range("B89").Value = "Invert Text With Gravity"
$ws.Range("C89").Value = "Invertir Texto Con La Gravedad"

$ws.Range("A90").Value = "EditText"
$ws.Range("B90").Value = "Edit Text"
$ws.Range("C90").Value = "Editar Texto"

# ---------------------------------------------------------------------
# 2. Fill in the "Collision" row that was inserted in step 0.
# ---------------------------------------------------------------------

$ws.Range("A63").Value = "Collision"
$ws.Range("B63").Value = "Collision"
$ws.Range("C63").Value = "Colisión"

# ---------------------------------------------------------------------
# 3. Restore the final selection state recorded in the saved file.
# ---------------------------------------------------------------------

$ws.Range("D63").Select()

# ---------------------------------------------------------------------
# 4. Page setup tweak captured by the authoring session.
# ---------------------------------------------------------------------

$ws.PageSetup.Orientation = 1
